# PSP 박영우 UPDATE 11.12
# Applies the recorded changes to the Time Recording Log sheet:
#  - Clears the now-unused "Date" column header (A5)
#  - Fixes the date on the existing "SAD 스켈레톤 코드..." entry (row 13)
#  - Replaces row 14 with a new "Entity Class 설계 및 산출물 작성" entry (11/5)
#  - Adds two brand-new log entries for 11/11 and 11/12 (rows 15 and 16)
#  - Normalizes the formatting/formulas of the following blank rows (17, 18)
#  - Moves the active selection to F20

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Column A no longer has a "Date" header; clear it and drop back to the
#    plain (unbordered) style used elsewhere for blank label cells.
$ws.Range("B3").Copy()
$ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A5").ClearContents()

# 2. Row 13: same activity, corrected date (9/25 -> 10/25).
$ws.Range("A13").Value = 43763

# 3. Row 16: new entry - ERD design (11/12). Its activity text is entered
#    first so it lands before "Entity Class..." in the shared-string table,
#    matching the order produced by the original edit.
$ws.Range("A16").Value = 43781
$ws.Range("B16").Value = 0.41666666666666669
$ws.Range("C16").Value = 0.5625
$ws.Range("D16").Value = 20
$ws.Range("E16").Formula = "=(C16-B16)*1440-D16"
$ws.Range("F16").Value = "ERD 설계 및 산출물 작성"

# 4. Row 14: new entry - Entity Class design (11/5), replacing old data.
$ws.Range("A14").Value = 43774
$ws.Range("B14").Value = 0.5625
$ws.Range("C14").Value = 0.60416666666666663
$ws.Range("D14").Value = 0
$ws.Range("E14").Formula = "=(C14-B14)*1440-D14"
$ws.Range("F14").Value = "Entity Class 설계 및 산출물 작성"

# 5. Row 15: new entry - Entity Class design continued (11/11).
$ws.Range("A15").Value = 43780
$ws.Range("B15").Value = 0.66666666666666663
$ws.Range("C15").Value = 0.71527777777777779
$ws.Range("D15").Value = 0
$ws.Range("E15").Formula = "=(C15-B15)*1440-D15"
$ws.Range("F15").Value = "Entity Class 설계 및 산출물 작성"

# 6. Rows 17-18: still-blank template rows now match the plain "no border"
#    style already used by the rows further below (19, 20, ...), with
#    explicit 0 interruption time and the delta-time formula filled in.
$ws.Range("F17").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D17").Value = 0
$ws.Range("E17").Formula = "=(C17-B17)*1440-D17"

$ws.Range("F18").Copy()
$ws.Range("A18:C18").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D18").Value = 0
$ws.Range("E18").Formula = "=(C18-B18)*1440-D18"

# 7. Move the active selection from G20 to F20.
$ws.Range("F20").Select()
